# Apply data updates to the 'Resumo Inscricoes' worksheet
# Updates Inscritos (E), Pagos (F) and Inscricoes homologadas (H) counts
# for several rows, per the upstream data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 113
$ws.Range("F2").Value = 81
$ws.Range("H2").Value = 87
$ws.Range("E6").Value = 52
$ws.Range("E10").Value = 703
$ws.Range("F10").Value = 393
$ws.Range("H10").Value = 488
$ws.Range("E12").Value = 702
$ws.Range("F12").Value = 423
$ws.Range("H12").Value = 509
$ws.Range("E13").Value = 167
$ws.Range("E16").Value = 233
$ws.Range("F22").Value = 109
$ws.Range("H22").Value = 151
$ws.Range("E25").Value = 332
$ws.Range("F25").Value = 184
$ws.Range("H25").Value = 244
$ws.Range("E26").Value = 199
$ws.Range("F26").Value = 123
$ws.Range("H26").Value = 148
$ws.Range("E27").Value = 384
$ws.Range("F27").Value = 211
$ws.Range("H27").Value = 293
$ws.Range("E28").Value = 229
$ws.Range("F28").Value = 113
$ws.Range("H28").Value = 165
$ws.Range("E32").Value = 218
$ws.Range("E33").Value = 331
$ws.Range("F33").Value = 180
$ws.Range("H33").Value = 270
$ws.Range("F34").Value = 177
$ws.Range("H34").Value = 215
$ws.Range("E35").Value = 185
$ws.Range("F35").Value = 127
$ws.Range("H35").Value = 154
$ws.Range("E39").Value = 201
$ws.Range("F39").Value = 104
$ws.Range("H39").Value = 155
$ws.Range("E41").Value = 445
$ws.Range("F41").Value = 223
$ws.Range("H41").Value = 315
$ws.Range("E42").Value = 466
$ws.Range("F42").Value = 264
$ws.Range("H42").Value = 325
$ws.Range("E43").Value = 144
$ws.Range("F43").Value = 82
$ws.Range("H43").Value = 109
$ws.Range("E44").Value = 373
$ws.Range("F44").Value = 195
$ws.Range("H44").Value = 263
$ws.Range("E47").Value = 540
$ws.Range("E48").Value = 272
$ws.Range("F48").Value = 132
$ws.Range("H48").Value = 176
$ws.Range("E49").Value = 340
$ws.Range("E50").Value = 285
$ws.Range("F51").Value = 133
$ws.Range("H51").Value = 207

$wb.Save()
